# PFPT_QTR_FIN.xlsx quarterly-financials update.
# Two new fiscal quarters (ending 2018-12-31 and 2018-09-30) are added as
# the new leftmost data columns D and E; the previously-existing quarter
# columns D:K slide right to F:M intact. A handful of trailing-quarter
# figures (columns H/I, the former F/G) were also restated in this pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the 2 new columns; Excel shifts D:K -> F:M, values/types intact.
$ws.Columns("D:E").Insert()

# The freshly inserted D:E columns inherit column C format; re-stamp them
# with the number/date format used by the rest of the data block (column F).
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- New quarter data (column D = latest quarter, column E = prior quarter) ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 198500
$ws.Range("E8").Value = 184200
$ws.Range("D9").Value = 53000
$ws.Range("E9").Value = 50900
$ws.Range("D10").Value = 145500
$ws.Range("E10").Value = 133300
$ws.Range("D12").Value = 48200
$ws.Range("E12").Value = 45900
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 3800
$ws.Range("E15").Value = 4000
$ws.Range("D17").Value = 219500
$ws.Range("E17").Value = 210700
$ws.Range("D18").Value = -21000
$ws.Range("E18").Value = -26500
$ws.Range("D20").Value = 500
$ws.Range("E20").Value = -400
$ws.Range("D21").Value = -1300
$ws.Range("E21").Value = -7500
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = 9100
$ws.Range("D23").Value = -20500
$ws.Range("E23").Value = -36100
$ws.Range("D24").Value = 700
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -21200
$ws.Range("E26").Value = -36100
$ws.Range("D27").Value = -21200
$ws.Range("E27").Value = -36100
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -500
$ws.Range("E32").Value = 400
$ws.Range("D33").Value = -21200
$ws.Range("E33").Value = -36100
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -21200
$ws.Range("E35").Value = -36100
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 185400
$ws.Range("E41").Value = 153000
$ws.Range("D42").Value = 46300
$ws.Range("E42").Value = 37400
$ws.Range("D43").Value = 199200
$ws.Range("E43").Value = 143900
$ws.Range("D44").Value = 500
$ws.Range("E44").Value = 400
$ws.Range("D45").Value = 56100
$ws.Range("E45").Value = 53900
$ws.Range("D46").Value = 487400
$ws.Range("E46").Value = 388500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 70600
$ws.Range("E48").Value = 73500
$ws.Range("D49").Value = 597100
$ws.Range("E49").Value = 607800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 77900
$ws.Range("E52").Value = 67800
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1233000
$ws.Range("E54").Value = 1137700
$ws.Range("D57").Value = 20200
$ws.Range("E57").Value = 14900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 584300
$ws.Range("E59").Value = 526400
$ws.Range("D60").Value = 604500
$ws.Range("E60").Value = 541400
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 116000
$ws.Range("E62").Value = 92400
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 720500
$ws.Range("E66").Value = 633900
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -595400
$ws.Range("E72").Value = -574200
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 512500
$ws.Range("E76").Value = 503800
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -21200
$ws.Range("E81").Value = -36100
$ws.Range("D83").Value = 19200
$ws.Range("E83").Value = 19400
$ws.Range("D89").Value = 55100
$ws.Range("E89").Value = 64700
$ws.Range("D91").Value = -6400
$ws.Range("E91").Value = -6500
$ws.Range("D94").Value = -22800
$ws.Range("E94").Value = -13700
$ws.Range("D100").Value = -7100
$ws.Range("E100").Value = -4600
$ws.Range("D101").Value = -400
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 24800
$ws.Range("E102").Value = 46300

# --- Restated figures for the two quarters now in columns H and I (and J for row 91) ---
$ws.Range("H8").Value = 146900
$ws.Range("I8").Value = 134700
$ws.Range("H10").Value = 106400
$ws.Range("I10").Value = 98700
$ws.Range("H17").Value = 158700
$ws.Range("I17").Value = 148300
$ws.Range("H18").Value = -11800
$ws.Range("I18").Value = -13600
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = -2600
$ws.Range("H23").Value = -19900
$ws.Range("I23").Value = -18500
$ws.Range("H26").Value = -8600
$ws.Range("I26").Value = -19500
$ws.Range("H27").Value = -8600
$ws.Range("I27").Value = -19500
$ws.Range("H33").Value = -6500
$ws.Range("I33").Value = -19500
$ws.Range("H35").Value = -6500
$ws.Range("I35").Value = -19500
$ws.Range("H81").Value = -6500
$ws.Range("I81").Value = -19500
$ws.Range("I91").Value = -11900
$ws.Range("J91").Value = -10600
